$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 287.875
$ws.Range("I2").Value = 287.875
$ws.Range("K2").Value = 287.875
$ws.Range("M2").Value = -174.875
$ws.Range("H4").Value = 110
$ws.Range("I4").Value = 80
$ws.Range("K4").Value = 80
$ws.Range("M4").Value = 34
$ws.Range("H17").Value = 1092.9459
$ws.Range("J17").Value = 1092.9459
$ws.Range("L17").Value = 3278.8377
$ws.Range("N17").Value = -3614.8377
$ws.Range("H33").Value = 586.7826
$ws.Range("I33").Value = 185.53334
$ws.Range("J33").Value = 1339.125
$ws.Range("K33").Value = 185.53334
$ws.Range("L33").Value = 1339.125
$ws.Range("M33").Value = 43.46665999999999
$ws.Range("N33").Value = -1797.125
$ws.Range("H51").Value = 3055.4443
$ws.Range("J51").Value = 3750
$ws.Range("L51").Value = 3750
$ws.Range("N51").Value = -4718
$ws.Range("H74").Value = 4074.625
$ws.Range("I74").Value = 4074.625
$ws.Range("K74").Value = 4074.625
$ws.Range("M74").Value = -3138.625
$ws.Range("H77").Value = 4074.625
$ws.Range("I77").Value = 4074.625
$ws.Range("K77").Value = 20373.125
$ws.Range("M77").Value = -15693.125
$ws.Range("H86").Value = 4622.5557
$ws.Range("I86").Value = 4800
$ws.Range("J86").Value = 4400.75
$ws.Range("K86").Value = 4800
$ws.Range("L86").Value = 4400.75
$ws.Range("M86").Value = -3677
$ws.Range("N86").Value = -6646.75
$ws.Range("H89").Value = 4622.5557
$ws.Range("I89").Value = 4800
$ws.Range("J89").Value = 4400.75
$ws.Range("K89").Value = 24000
$ws.Range("L89").Value = 22003.75
$ws.Range("M89").Value = -18384
$ws.Range("N89").Value = -33235.75
$ws.Range("H92").Value = 7422
$ws.Range("I92").Value = 3186.2856
$ws.Range("K92").Value = 3186.2856
$ws.Range("M92").Value = -1938.2856
$ws.Range("H97").Value = 6113.1113
$ws.Range("J97").Value = 6113.1113
$ws.Range("L97").Value = 18339.3339
$ws.Range("N97").Value = -19331.3339
$ws.Range("H99").Value = 936.625
$ws.Range("I99").Value = 356.14285
$ws.Range("J99").Value = 5000
$ws.Range("K99").Value = 1068.42855
$ws.Range("L99").Value = 15000
$ws.Range("M99").Value = 429.5714499999999
$ws.Range("N99").Value = -17996
$ws.Range("H107").Value = 307.3846
$ws.Range("I107").Value = 307.3846
$ws.Range("K107").Value = 307.3846
$ws.Range("M107").Value = 1612.6154
$ws.Range("H108").Value = 52000
$ws.Range("J108").Value = 52000
$ws.Range("L108").Value = 52000
$ws.Range("N108").Value = -59680
$ws.Range("H112").Value = 2560.875
$ws.Range("I112").Value = 1533.3334
$ws.Range("J112").Value = 3177.4
$ws.Range("K112").Value = 4600.0002
$ws.Range("L112").Value = 9532.200000000001
$ws.Range("M112").Value = -3492.0002
$ws.Range("N112").Value = -11748.2
$ws.Range("H113").Value = 2188
$ws.Range("I113").Value = 2188
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2188
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1066
$ws.Range("N113").ClearContents()
$ws.Range("H129").Value = 1745.6923
$ws.Range("I129").Value = 854.8889
$ws.Range("K129").Value = 2564.6667
$ws.Range("M129").Value = 2435.3333
$ws.Range("H131").Value = 6611.5864
$ws.Range("I131").Value = 5336.9
$ws.Range("J131").Value = 9444.223
$ws.Range("K131").Value = 16010.7
$ws.Range("L131").Value = 28332.669
$ws.Range("M131").Value = -10970.7
$ws.Range("N131").Value = -38412.669
$ws.Range("H132").Value = 3843.8
$ws.Range("I132").Value = 3242.25
$ws.Range("J132").Value = 6250
$ws.Range("K132").Value = 9726.75
$ws.Range("L132").Value = 18750
$ws.Range("M132").Value = -7196.75
$ws.Range("N132").Value = -23810
$ws.Range("H137").Value = 9277
$ws.Range("I137").Value = 9606
$ws.Range("K137").Value = 28818
$ws.Range("M137").Value = -26268
$ws.Range("H138").Value = 2243.519
$ws.Range("I138").Value = 1168.8096
$ws.Range("J138").Value = 2632.638
$ws.Range("K138").Value = 3506.4288
$ws.Range("L138").Value = 7897.914
$ws.Range("M138").Value = 1633.5712
$ws.Range("N138").Value = -18177.914
$ws.Range("H141").Value = 6905.08
$ws.Range("I141").Value = 4509.846
$ws.Range("J141").Value = 9499.916999999999
$ws.Range("K141").Value = 13529.538
$ws.Range("L141").Value = 28499.751
$ws.Range("M141").Value = -8349.537999999999
$ws.Range("N141").Value = -38859.751

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2926.6316
$ws.Range("I2").Value = 1061.6
$ws.Range("J2").Value = 4998.8887
$ws.Range("K2").Value = 1061.6
$ws.Range("L2").Value = 4998.8887
$ws.Range("M2").Value = -948.5999999999999
$ws.Range("N2").Value = -5224.8887
$ws.Range("H32").Value = 3476.9714
$ws.Range("I32").Value = 3476.9714
$ws.Range("K32").Value = 3476.9714
$ws.Range("M32").Value = -3189.9714
$ws.Range("H33").Value = 26750.25
$ws.Range("J33").Value = 30000
$ws.Range("L33").Value = 30000
$ws.Range("N33").Value = -30658
$ws.Range("H45").Value = 724
$ws.Range("I45").Value = 724
$ws.Range("K45").Value = 724
$ws.Range("M45").Value = -347
$ws.Range("H74").Value = 1970.4736
$ws.Range("I74").Value = 2258.4
$ws.Range("J74").Value = 890.75
$ws.Range("K74").Value = 2258.4
$ws.Range("L74").Value = 890.75
$ws.Range("M74").Value = -1384.4
$ws.Range("N74").Value = -2638.75
$ws.Range("H77").Value = 1970.4736
$ws.Range("I77").Value = 2258.4
$ws.Range("J77").Value = 890.75
$ws.Range("K77").Value = 11292
$ws.Range("L77").Value = 4453.75
$ws.Range("M77").Value = -6924
$ws.Range("N77").Value = -13189.75
$ws.Range("H88").Value = 1722.6666
$ws.Range("I88").Value = 635
$ws.Range("K88").Value = 635
$ws.Range("M88").Value = -229
$ws.Range("H91").Value = 1722.6666
$ws.Range("I91").Value = 635
$ws.Range("K91").Value = 635
$ws.Range("M91").Value = 769
$ws.Range("H102").Value = 3196
$ws.Range("I102").Value = 2241.2307
$ws.Range("J102").Value = 7333.3335
$ws.Range("K102").Value = 2241.2307
$ws.Range("L102").Value = 7333.3335
$ws.Range("M102").Value = -619.2307000000001
$ws.Range("N102").Value = -10577.3335
$ws.Range("H110").Value = 1023.9286
$ws.Range("I110").Value = 1023.9286
$ws.Range("K110").Value = 1023.9286
$ws.Range("M110").Value = 1021.0714
$ws.Range("H116").Value = 2926.6316
$ws.Range("I116").Value = 1061.6
$ws.Range("J116").Value = 4998.8887
$ws.Range("K116").Value = 1061.6
$ws.Range("L116").Value = 4998.8887
$ws.Range("M116").Value = 1232.4
$ws.Range("N116").Value = -9586.8887
$ws.Range("H122").Value = 1435.44
$ws.Range("I122").Value = 1247.4762
$ws.Range("J122").Value = 2422.25
$ws.Range("K122").Value = 3742.4286
$ws.Range("L122").Value = 7266.75
$ws.Range("M122").Value = -1292.4286
$ws.Range("N122").Value = -12166.75
$ws.Range("H130").Value = 52463.168
$ws.Range("I130").Value = 53695
$ws.Range("J130").Value = 49999.5
$ws.Range("K130").Value = 53695
$ws.Range("L130").Value = 49999.5
$ws.Range("M130").Value = -48675
$ws.Range("N130").Value = -60039.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2926.6316
$ws.Range("I3").Value = 1061.6
$ws.Range("J3").Value = 4998.8887
$ws.Range("K3").Value = 1061.6
$ws.Range("L3").Value = 4998.8887
$ws.Range("M3").Value = -947.5999999999999
$ws.Range("N3").Value = -5226.8887
$ws.Range("H31").Value = 1500
$ws.Range("I31").Value = 1500
$ws.Range("K31").Value = 1500
$ws.Range("M31").Value = -1248
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()
$ws.Range("H37").Value = 450.5
$ws.Range("I37").Value = 450.5
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 450.5
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -313.5
$ws.Range("N37").ClearContents()
$ws.Range("H86").Value = 4600
$ws.Range("I86").Value = 4500
$ws.Range("J86").Value = 4800
$ws.Range("K86").Value = 4500
$ws.Range("L86").Value = 4800
$ws.Range("M86").Value = -3377
$ws.Range("N86").Value = -7046
$ws.Range("H89").Value = 4600
$ws.Range("I89").Value = 4500
$ws.Range("J89").Value = 4800
$ws.Range("K89").Value = 22500
$ws.Range("L89").Value = 24000
$ws.Range("M89").Value = -16884
$ws.Range("N89").Value = -35232
$ws.Range("H94").Value = 2200.3794
$ws.Range("I94").Value = 1967.7142
$ws.Range("J94").Value = 2811.125
$ws.Range("K94").Value = 1967.7142
$ws.Range("L94").Value = 2811.125
$ws.Range("M94").Value = -1516.7142
$ws.Range("N94").Value = -3713.125
$ws.Range("H99").Value = 4190.5625
$ws.Range("I99").Value = 2948.6667
$ws.Range("J99").Value = 5787.2856
$ws.Range("K99").Value = 2948.6667
$ws.Range("L99").Value = 5787.2856
$ws.Range("M99").Value = -1450.6667
$ws.Range("N99").Value = -8783.285599999999
$ws.Range("H107").Value = 4006.0571
$ws.Range("I107").Value = 1486.05
$ws.Range("J107").Value = 7366.067
$ws.Range("K107").Value = 1486.05
$ws.Range("L107").Value = 7366.067
$ws.Range("M107").Value = 433.95
$ws.Range("N107").Value = -11206.067
$ws.Range("H111").Value = 62500
$ws.Range("J111").Value = 62500
$ws.Range("L111").Value = 62500
$ws.Range("N111").Value = -70680

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 109.9
$ws.Range("I7").Value = 71.09524
$ws.Range("J7").Value = 200.44444
$ws.Range("K7").Value = 71.09524
$ws.Range("L7").Value = 200.44444
$ws.Range("M7").Value = 41.90476
$ws.Range("N7").Value = -426.44444
$ws.Range("H16").Value = 1804.75
$ws.Range("I16").Value = 1804.75
$ws.Range("K16").Value = 1804.75
$ws.Range("M16").Value = -1517.75
$ws.Range("H22").Value = 1197.069
$ws.Range("I22").Value = 903.7059
$ws.Range("J22").Value = 1612.6666
$ws.Range("K22").Value = 903.7059
$ws.Range("L22").Value = 1612.6666
$ws.Range("M22").Value = -553.7059
$ws.Range("N22").Value = -2312.6666
$ws.Range("H31").Value = 3830.44
$ws.Range("I31").Value = 1593.5238
$ws.Range("K31").Value = 1593.5238
$ws.Range("M31").Value = -1298.5238
$ws.Range("H34").Value = 3830.44
$ws.Range("I34").Value = 1593.5238
$ws.Range("K34").Value = 1593.5238
$ws.Range("M34").Value = -1391.5238
$ws.Range("H99").Value = 2512.5
$ws.Range("I99").Value = 2015
$ws.Range("K99").Value = 2015
$ws.Range("M99").Value = -517
$ws.Range("H107").Value = 35944.17
$ws.Range("I107").Value = 167081.83
$ws.Range("J107").Value = 1734.3478
$ws.Range("K107").Value = 167081.83
$ws.Range("L107").Value = 1734.3478
$ws.Range("M107").Value = -165161.83
$ws.Range("N107").Value = -5574.3478
$ws.Range("H113").Value = 1804.75
$ws.Range("I113").Value = 1804.75
$ws.Range("K113").Value = 1804.75
$ws.Range("M113").Value = 365.25
$ws.Range("H118").Value = 74991.664
$ws.Range("J118").Value = 74991.664
$ws.Range("L118").Value = 74991.664
$ws.Range("N118").Value = -78305.664
$ws.Range("H126").Value = 2512.5
$ws.Range("I126").Value = 2015
$ws.Range("K126").Value = 6045
$ws.Range("M126").Value = -3575
$ws.Range("H134").Value = 2071.484
$ws.Range("I134").Value = 1949.1724
$ws.Range("K134").Value = 5847.5172
$ws.Range("M134").Value = -3312.5172
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 643.4
$ws.Range("J38").Value = 907.7143
$ws.Range("L38").Value = 2723.1429
$ws.Range("N38").Value = -3417.1429
$ws.Range("H55").Value = 1000373.2
$ws.Range("I55").Value = 1000373.2
$ws.Range("K55").Value = 3001119.6
$ws.Range("M55").Value = -3000942.6
$ws.Range("H86").Value = 500
$ws.Range("J86").Value = 500
$ws.Range("L86").Value = 1500
$ws.Range("N86").Value = -3872
$ws.Range("H89").Value = 500
$ws.Range("J89").Value = 500
$ws.Range("L89").Value = 4500
$ws.Range("N89").Value = -16356
$ws.Range("H112").Value = 203901.6
$ws.Range("I112").Value = 335902.34
$ws.Range("K112").Value = 1007707.02
$ws.Range("M112").Value = -1006599.02
$ws.Range("H113").Value = 1110.7646
$ws.Range("J113").Value = 1195.2142
$ws.Range("L113").Value = 3585.6426
$ws.Range("N113").Value = -7925.642599999999
$ws.Range("H121").Value = 72079.92999999999
$ws.Range("J121").Value = 799.75
$ws.Range("L121").Value = 2399.25
$ws.Range("N121").Value = -5019.25
$ws.Range("H133").Value = 5651
$ws.Range("I133").Value = 3477
$ws.Range("J133").Value = 9999
$ws.Range("K133").Value = 10431
$ws.Range("L133").Value = 29997
$ws.Range("M133").Value = -5371
$ws.Range("N133").Value = -40117
$ws.Range("H140").Value = 85076.25
$ws.Range("I140").Value = 101303.5
$ws.Range("J140").Value = 3940
$ws.Range("K140").Value = 303910.5
$ws.Range("L140").Value = 11820
$ws.Range("M140").Value = -298730.5
$ws.Range("N140").Value = -22180
$ws.Range("H141").Value = 1853.1111
$ws.Range("I141").Value = 1853.1111
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 5559.3333
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -379.3333000000002
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H55").Value = 26277.334
$ws.Range("I55").Value = 16166.667
$ws.Range("J55").Value = 31332.666
$ws.Range("K55").Value = 16166.667
$ws.Range("L55").Value = 31332.666
$ws.Range("M55").Value = -15839.667
$ws.Range("N55").Value = -31986.666
$ws.Range("H70").Value = 9844.507
$ws.Range("I70").Value = 6135.4
$ws.Range("K70").Value = 6135.4
$ws.Range("M70").Value = -5865.4
$ws.Range("H73").Value = 9844.507
$ws.Range("I73").Value = 6135.4
$ws.Range("K73").Value = 6135.4
$ws.Range("M73").Value = -5199.4
$ws.Range("H80").Value = 23167.334
$ws.Range("I80").Value = 55002.5
$ws.Range("J80").Value = 7249.75
$ws.Range("K80").Value = 55002.5
$ws.Range("L80").Value = 7249.75
$ws.Range("M80").Value = -54004.5
$ws.Range("N80").Value = -9245.75
$ws.Range("H83").Value = 23167.334
$ws.Range("I83").Value = 55002.5
$ws.Range("J83").Value = 7249.75
$ws.Range("K83").Value = 275012.5
$ws.Range("L83").Value = 36248.75
$ws.Range("M83").Value = -270020.5
$ws.Range("N83").Value = -46232.75
$ws.Range("H102").Value = 1699.9333
$ws.Range("I102").Value = 1699.9333
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1699.9333
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -77.93329999999992
$ws.Range("N102").ClearContents()
$ws.Range("H113").Value = 4726.727
$ws.Range("I113").Value = 1341.1212
$ws.Range("J113").Value = 9805.137000000001
$ws.Range("K113").Value = 1341.1212
$ws.Range("L113").Value = 9805.137000000001
$ws.Range("M113").Value = 828.8788
$ws.Range("N113").Value = -14145.137
$ws.Range("H117").Value = 49499.9
$ws.Range("J117").Value = 49499.9
$ws.Range("L117").Value = 49499.9
$ws.Range("N117").Value = -56383.9
$ws.Range("H122").Value = 11355.429
$ws.Range("I122").Value = 11459.77
$ws.Range("K122").Value = 34379.31
$ws.Range("M122").Value = -31929.31
$ws.Range("H126").Value = 3091.6316
$ws.Range("I126").Value = 2938.8572
$ws.Range("J126").Value = 3519.4
$ws.Range("K126").Value = 8816.571599999999
$ws.Range("L126").Value = 10558.2
$ws.Range("M126").Value = -6346.571599999999
$ws.Range("N126").Value = -15498.2
$ws.Range("H132").Value = 2611.4849
$ws.Range("I132").Value = 2627.5334
$ws.Range("K132").Value = 7882.600199999999
$ws.Range("M132").Value = -5352.600199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 1349.1666
$ws.Range("I9").Value = 1024
$ws.Range("K9").Value = 1024
$ws.Range("M9").Value = -800
$ws.Range("H16").Value = 494.46155
$ws.Range("I16").Value = 514.2222
$ws.Range("K16").Value = 514.2222
$ws.Range("M16").Value = -344.2222
$ws.Range("H32").Value = 9833
$ws.Range("I32").Value = 9833
$ws.Range("K32").Value = 9833
$ws.Range("M32").Value = -9516
$ws.Range("H40").Value = 5166.3335
$ws.Range("I40").Value = 3539.0557
$ws.Range("J40").Value = 7119.067
$ws.Range("K40").Value = 3539.0557
$ws.Range("L40").Value = 7119.067
$ws.Range("M40").Value = -3403.0557
$ws.Range("N40").Value = -7391.067
$ws.Range("H55").Value = 1205.15
$ws.Range("I55").Value = 275.7857
$ws.Range("K55").Value = 275.7857
$ws.Range("M55").Value = -102.7857
$ws.Range("H61").Value = 4488.316
$ws.Range("I61").Value = 1137.8
$ws.Range("K61").Value = 1137.8
$ws.Range("M61").Value = -935.8
$ws.Range("H93").Value = 4418.5
$ws.Range("I93").Value = 1796.5
$ws.Range("J93").Value = 6166.5
$ws.Range("K93").Value = 1796.5
$ws.Range("L93").Value = 6166.5
$ws.Range("M93").Value = -548.5
$ws.Range("N93").Value = -8662.5
$ws.Range("H100").Value = 4910.1113
$ws.Range("I100").Value = 2365.1667
$ws.Range("K100").Value = 2365.1667
$ws.Range("M100").Value = -1824.1667
$ws.Range("H113").Value = 4488.316
$ws.Range("I113").Value = 1137.8
$ws.Range("K113").Value = 1137.8
$ws.Range("M113").Value = 1032.2
$ws.Range("H118").Value = 210999.5
$ws.Range("J118").Value = 210999.5
$ws.Range("L118").Value = 210999.5
$ws.Range("N118").Value = -214313.5
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()
$ws.Range("H122").Value = 3776
$ws.Range("I122").Value = 3359.5789
$ws.Range("J122").Value = 5358.4
$ws.Range("K122").Value = 10078.7367
$ws.Range("L122").Value = 16075.2
$ws.Range("M122").Value = -7628.736699999999
$ws.Range("N122").Value = -20975.2
$ws.Range("H132").Value = 3905.907
$ws.Range("I132").Value = 4026.55
$ws.Range("J132").Value = 3801
$ws.Range("K132").Value = 12079.65
$ws.Range("L132").Value = 11403
$ws.Range("M132").Value = -9549.650000000001
$ws.Range("N132").Value = -16463

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 40122000
$ws.Range("J5").Value = 66688332
$ws.Range("L5").Value = 66688332
$ws.Range("N5").Value = -66688556
$ws.Range("H24").Value = 14997.5
$ws.Range("I24").Value = 14997.5
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 14997.5
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = -14767.5
$ws.Range("N24").ClearContents()
$ws.Range("H62").Value = 11857.429
$ws.Range("I62").Value = 11002
$ws.Range("K62").Value = 11002
$ws.Range("M62").Value = -10378
$ws.Range("H65").Value = 11857.429
$ws.Range("I65").Value = 11002
$ws.Range("K65").Value = 55010
$ws.Range("M65").Value = -51890
$ws.Range("H70").Value = 45832
$ws.Range("I70").Value = 30000
$ws.Range("J70").Value = 48998.4
$ws.Range("K70").Value = 30000
$ws.Range("L70").Value = 48998.4
$ws.Range("M70").Value = -29685
$ws.Range("N70").Value = -49628.4
$ws.Range("H73").Value = 45832
$ws.Range("I73").Value = 30000
$ws.Range("J73").Value = 48998.4
$ws.Range("K73").Value = 30000
$ws.Range("L73").Value = 48998.4
$ws.Range("M73").Value = -28908
$ws.Range("N73").Value = -51182.4
$ws.Range("H96").Value = 3419.2222
$ws.Range("I96").Value = 2267
$ws.Range("J96").Value = 3995.3333
$ws.Range("K96").Value = 2267
$ws.Range("L96").Value = 3995.3333
$ws.Range("M96").Value = -894
$ws.Range("N96").Value = -6741.3333
$ws.Range("H100").Value = 411.66666
$ws.Range("I100").Value = 459.25
$ws.Range("J100").Value = 373.6
$ws.Range("K100").Value = 918.5
$ws.Range("L100").Value = 747.2
$ws.Range("M100").Value = -377.5
$ws.Range("N100").Value = -1829.2
$ws.Range("H107").Value = 4102.6665
$ws.Range("I107").Value = 1930.7894
$ws.Range("J107").Value = 12355.8
$ws.Range("K107").Value = 5792.3682
$ws.Range("L107").Value = 37067.39999999999
$ws.Range("M107").Value = -3872.3682
$ws.Range("N107").Value = -40907.39999999999
$ws.Range("H112").Value = 40999.8
$ws.Range("J112").Value = 40999.8
$ws.Range("L112").Value = 40999.8
$ws.Range("N112").Value = -43953.8
$ws.Range("H113").Value = 1431.5385
$ws.Range("I113").Value = 1429
$ws.Range("J113").Value = 1440
$ws.Range("K113").Value = 4287
$ws.Range("L113").Value = 4320
$ws.Range("M113").Value = -2117
$ws.Range("N113").Value = -8660
$ws.Range("H126").Value = 1441.7142
$ws.Range("I126").Value = 1199.25
$ws.Range("J126").Value = 1765
$ws.Range("K126").Value = 3597.75
$ws.Range("L126").Value = 5295
$ws.Range("M126").Value = -1127.75
$ws.Range("N126").Value = -10235
